# Limpieza del dataframe de productos:
# La columna "categoria" (C) tenía valores alternados inconsistentes dentro
# de bloques que en realidad pertenecen a una sola categoria. Se corrige
# para que cada bloque de filas tenga una categoria uniforme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bloques de filas (1-based, incluyendo encabezado en la fila 1) y la
# categoria correcta para cada bloque, según el dataset limpio final.
$blocks = @(
    @{ Start = 2;  End = 52;  Categoria = "Alimentos" },
    @{ Start = 53; End = 58;  Categoria = "Limpieza" },
    @{ Start = 59; End = 90;  Categoria = "Alimentos" },
    @{ Start = 91; End = 101; Categoria = "Limpieza" }
)

foreach ($block in $blocks) {
    for ($row = $block.Start; $row -le $block.End; $row++) {
        $cell = $ws.Cells.Item($row, 3)
        if ($cell.Value -ne $block.Categoria) {
            $cell.Value = $block.Categoria
        }
    }
}
